# Apply weekly price-record update for rows 97-199 (Hortaliza / Ciboulette sheet).
# Row 97 (old outlier date 44596, unit "$/atado") is dropped; rows 98-199 shift up by
# one position, and a new weekly record is appended as the final row (199).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> [D(Fecha), I(Calidad), J(Volumen), K(Precio minimo), L(Precio maximo),
#               M(Precio promedio ponderado), N(Unidad de comercializacion), O(Origen),
#               P(Precio $/Kg), Q(Kg o Unidades)]
$data = @{
    97 = @(44740, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    98 = @(44245, "Primera", 50, 7000, 7000, 7000, "`$/docena de atados", "Provincia de Cautín", 2333, 3)
    99 = @(44397, "Primera", 240, 4000, 4000, 4000, "`$/docena de atados", "Provincia de Chacabuco", 1333, 3)
    100 = @(44698, "Primera", 240, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    101 = @(44522, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    102 = @(44539, "Primera", 120, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    103 = @(44495, "Primera", 240, 2500, 3000, 2750, "`$/docena de atados", "Región Metropolitana", 917, 3)
    104 = @(44747, "Primera", 240, 2500, 3000, 2750, "`$/docena de atados", "Región Metropolitana", 917, 3)
    105 = @(44971, "Primera", 240, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    106 = @(44442, "Primera", 250, 4500, 4500, 4500, "`$/docena de atados", "Región Metropolitana", 1500, 3)
    107 = @(44293, "Primera", 40, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    108 = @(44253, "Primera", 250, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    109 = @(44537, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    110 = @(44784, "Primera", 80, 4000, 4000, 4000, "`$/docena de atados", "Región Metropolitana", 1333, 3)
    111 = @(44719, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    112 = @(44974, "Primera", 240, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    113 = @(44879, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    114 = @(44476, "Primera", 120, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    115 = @(44326, "Primera", 40, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    116 = @(44998, "Primera", 40, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    117 = @(45049, "Primera", 40, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    118 = @(44775, "Primera", 240, 4000, 4000, 4000, "`$/docena de atados", "Región Metropolitana", 1333, 3)
    119 = @(44232, "Primera", 200, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    120 = @(44351, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    121 = @(44876, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    122 = @(44322, "Primera", 40, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    123 = @(44435, "Primera", 560, 4500, 5000, 4786, "`$/docena de atados", "Región Metropolitana", 1595, 3)
    124 = @(45009, "Primera", 240, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    125 = @(44609, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    126 = @(44540, "Primera", 280, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    127 = @(44407, "Primera", 270, 4500, 5000, 4778, "`$/docena de atados", "Región Metropolitana", 1593, 3)
    128 = @(44379, "Primera", 240, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    129 = @(45008, "Primera", 80, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    130 = @(44281, "Primera", 40, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667, 3)
    131 = @(44257, "Primera", 200, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    132 = @(44917, "Primera", 120, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    133 = @(44946, "Primera", 240, 3000, 3500, 3250, "`$/docena de atados", "Región Metropolitana", 1083, 3)
    134 = @(44699, "Primera", 40, 7000, 7000, 7000, "`$/docena de atados", "Provincia de Cautín", 2333, 3)
    135 = @(44357, "Primera", 80, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    136 = @(44481, "Primera", 300, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    137 = @(44897, "Primera", 220, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    138 = @(44838, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    139 = @(44964, "Primera", 240, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    140 = @(44376, "Primera", 240, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    141 = @(44771, "Primera", 240, 4000, 4000, 4000, "`$/docena de atados", "Región Metropolitana", 1333, 3)
    142 = @(44875, "Primera", 120, 2000, 2500, 2250, "`$/docena de atados", "Región Metropolitana", 750, 3)
    143 = @(44777, "Primera", 120, 4000, 4000, 4000, "`$/docena de atados", "Región Metropolitana", 1333, 3)
    144 = @(44159, "Primera", 250, 2500, 3000, 2760, "`$/docena de atados", "Región Metropolitana", 920, 3)
    145 = @(44966, "Primera", 80, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    146 = @(44761, "Primera", 240, 4000, 4500, 4250, "`$/docena de atados", "Región Metropolitana", 1417, 3)
    147 = @(44243, "Primera", 240, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667, 3)
    148 = @(44721, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    149 = @(44943, "Primera", 240, 3000, 3500, 3250, "`$/docena de atados", "Región Metropolitana", 1083, 3)
    150 = @(44411, "Primera", 250, 4500, 4500, 4500, "`$/docena de atados", "Región Metropolitana", 1500, 3)
    151 = @(44432, "Primera", 240, 5000, 5000, 5000, "`$/docena de atados", "Región Metropolitana", 1667, 3)
    152 = @(44746, "Primera", 40, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    153 = @(44992, "Primera", 240, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    154 = @(44364, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    155 = @(44418, "Primera", 240, 5000, 5000, 5000, "`$/docena de atados", "Región Metropolitana", 1667, 3)
    156 = @(44466, "Primera", 90, 4000, 4000, 4000, "`$/docena de atados", "Región Metropolitana", 1333, 3)
    157 = @(44845, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    158 = @(44434, "Primera", 80, 5000, 5000, 5000, "`$/docena de atados", "Región Metropolitana", 1667, 3)
    159 = @(44736, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    160 = @(44161, "Primera", 80, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    161 = @(44973, "Primera", 80, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    162 = @(44169, "Primera", 250, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    163 = @(44533, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    164 = @(44358, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    165 = @(44365, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    166 = @(44312, "Primera", 40, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    167 = @(44238, "Primera", 50, 5000, 5000, 5000, "`$/docena de atados", "Provincia de Cautín", 1667, 3)
    168 = @(44224, "Primera", 80, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    169 = @(44371, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    170 = @(44980, "Primera", 100, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    171 = @(44623, "Primera", 120, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    172 = @(44523, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    173 = @(45029, "Primera", 80, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    174 = @(44525, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    175 = @(44448, "Primera", 120, 4000, 4000, 4000, "`$/docena de atados", "Región Metropolitana", 1333, 3)
    176 = @(44390, "Primera", 240, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    177 = @(45055, "Primera", 240, 3000, 3500, 3250, "`$/docena de atados", "Región Metropolitana", 1083, 3)
    178 = @(44711, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    179 = @(44333, "Primera", 40, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    180 = @(44858, "Primera", 40, 2000, 2000, 2000, "`$/docena de atados", "Región Metropolitana", 667, 3)
    181 = @(44858, "Segunda", 40, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    182 = @(44348, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    183 = @(44477, "Primera", 240, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    184 = @(44306, "Primera", 150, 5500, 5500, 5500, "`$/docena de atados", "Provincia de Cautín", 1833, 3)
    185 = @(44553, "Primera", 120, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    186 = @(44529, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    187 = @(44386, "Primera", 240, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    188 = @(44162, "Primera", 250, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    189 = @(44599, "Primera", 80, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    190 = @(44663, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    191 = @(45002, "Primera", 240, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    192 = @(44530, "Primera", 320, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    193 = @(44614, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    194 = @(44978, "Primera", 150, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    195 = @(44638, "Primera", 120, 6000, 6000, 6000, "`$/docena de atados", "Provincia de Cautín", 2000, 3)
    196 = @(44708, "Primera", 240, 2500, 2500, 2500, "`$/docena de atados", "Región Metropolitana", 833, 3)
    197 = @(44166, "Primera", 250, 3000, 3000, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
    198 = @(45033, "Primera", 80, 3500, 3500, 3500, "`$/docena de atados", "Región Metropolitana", 1167, 3)
    199 = @(45062, "Primera", 240, 2500, 3500, 3000, "`$/docena de atados", "Región Metropolitana", 1000, 3)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 9).Value  = $vals[1]   # I - Calidad
    $ws.Cells.Item($r, 10).Value = $vals[2]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[3]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[4]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $vals[6]   # N - Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $vals[7]   # O - Origen
    $ws.Cells.Item($r, 16).Value = $vals[8]   # P - Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $vals[9]   # Q - Kg o Unidades
}
